# Update the "DigiKey / Rectangular Connectors" line item (row 11) to the
# new "BuyaPi / 40 Pin GPIO Connector Header" line item, refresh the
# shipping/import-fee values below it, update the subtotal, and update the
# grand-total cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: swap the DigiKey connector line for the BuyaPi header line ---
$ws.Range("A11").Value = "40 Pin GPIO Connector Header"
$ws.Range("B11").Value = "BuyaPi"
$ws.Range("C11").Value = "SKU: 412"
$ws.Range("E11").Value = 1.95
$ws.Range("H11").Value = "https://www.buyapi.ca/product/40-pin-gpio-connector-header/"

# --- Row 12 (Tax and Customs) / Row 13 (Shipping and Handling) updates ---
$ws.Range("E12").Value = 1.68
$ws.Range("E13").Value = 11

# --- Row 14: updated subtotal for this section ---
$ws.Range("G14").Value = 14.63

# --- Row 22: updated grand total ---
$ws.Range("F22").Value = 88.13

# --- Fix up the H11 hyperlink so it points at the new BuyaPi product page.
# Hyperlinks.Add() appends rather than replacing in-place, and this engine's
# Hyperlinks.Delete() clears the whole sheet collection, so clear it once and
# re-add all three links (H5, H6 unchanged; H11 updated). ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H5"), "https://www.sparkfun.com/products/14193?_ga=2.97662492.2095878335.1537831851-1448490607.1536633859")
$ws.Hyperlinks.Add($ws.Range("H6"), "https://www.sparkfun.com/products/116?_ga=2.161208638.2095878335.1537831851-1448490607.1536633859")
$ws.Hyperlinks.Add($ws.Range("H11"), "https://www.buyapi.ca/product/40-pin-gpio-connector-header/")

# --- Update the saved selection state to match the edited workbook ---
$ws.Range("H15").Select()
